# Update the status tracker worksheet:
# Row 4  (Task 3 - Home Page): Status -> "Not Started", Progress (%) -> 0
# Row 22 (Task 21 - Service List Page): Status -> "Not Started", Progress (%) -> 0
# Also move the active selection to G22 to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Android_UserApp_Status_Tracker")

$ws.Range("F4").Value = "Not Started"
$ws.Range("G4").Value = 0

$ws.Range("F22").Value = "Not Started"
$ws.Range("G22").Value = 0

$ws.Range("G22").Select()
